$wb = $excel.ActiveWorkbook

# --- Users sheet: update the logged-in user name ---
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Indrajeet Singh"

# --- Company sheet: update the default/test company name ---
$wsCompany = $wb.Worksheets.Item("Company")
$wsCompany.Range("A2").Value = "CapProviderTestCompany"
$wsCompany.Range("G3").Select() | Out-Null

# --- ActivityColumns sheet: rework the dashboard column list ---
$wsCols = $wb.Worksheets.Item("ActivityColumns")
$wsCols.Range("A3").Value = "Open"
$wsCols.Range("A4").Value = "Edit"
$wsCols.Range("A5").Value = "Date"
$wsCols.Range("A6").Value = "Company Name"
$wsCols.Range("A7").Value = "Type"
$wsCols.Range("A8").Value = "Tier"
$wsCols.Range("A9").Value = "Event/Task Type"
$wsCols.Range("A10").Value = "HL Contact"
$wsCols.Range("A11").Value = "Subject"
$wsCols.Range("A12").Value = "Companies Discussed"
$wsCols.Range("A13").Value = "Meeting/Call Notes"
$wsCols.Range("A14").Value = "External Contact"
$wsCols.Range("A15").Value = "External Contact Company"
$wsCols.Range("A15").Select() | Out-Null

# --- Users sheet stays the active/selected sheet & cell ---
$wsUsers.Activate() | Out-Null
$wsUsers.Range("A2").Select() | Out-Null
